# Commit: "Add files via upload"
# The Test-Cases sheet has an "Approved/Rejected" column (I) that already
# carries the value "Approved" for rows 25-52 (styled as text, numFmt "@").
# This edit fills in the same "Approved" value (with the same text style)
# for the remaining rows 2-24, and updates the saved view state (which
# cell range is selected / scrolled into view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 24 of column I ("Approved/Rejected") get the value
# "Approved", matching the style already used in I25:I52 (text number
# format so a future "Rejected" typed in stays literal text).
$rng = $ws.Range("I2:I24")
$rng.NumberFormat = "@"
$rng.Value = "Approved"

# Restore/scroll the window so the newly completed range is what's
# visible & selected (matches the sheetView/selection recorded after the
# edit: topLeftCell G5, activeCell I2, selection I2:I24).
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 7
$ws.Range("I2:I24").Select()
